$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (cell, new value) updates
# Generated from the commit diff: "Add data for 2022-06-01"
# (updates the running 2022 totals, column I, for June 2022 data)
$updates = @{}

$updates['Citywide Totals'] = @(
    @('I2', 2684),
    @('H3', 8347),
    @('I3', 2811),
    @('I4', 688),
    @('I5', 249),
    @('I6', 3190),
    @('H7', 25972),
    @('I7', 9622)
)

$updates['Uptown'] = @(
    @('I2', 29),
    @('I7', 104)
)

$updates['West Ridge'] = @(
    @('I3', 32),
    @('I7', 117)
)

$updates['Bridgeport'] = @(
    @('I3', 13),
    @('I7', 48)
)

$updates['Woodlawn'] = @(
    @('I5', 4),
    @('I7', 174)
)

$updates['North Lawndale'] = @(
    @('I3', 128),
    @('I6', 128),
    @('I7', 366)
)

$updates['South Deering'] = @(
    @('I3', 25),
    @('I7', 76)
)

$updates['New City'] = @(
    @('I2', 67),
    @('I3', 56)
)

$updates['By Neighborhood'] = @(
    @('I2', 93),
    @('I4', 39),
    @('I7', 320),
    @('I8', 614),
    @('I11', 160),
    @('I14', 48),
    @('I15', 125),
    @('I20', 241),
    @('I21', 57),
    @('I29', 642),
    @('I33', 452),
    @('I41', 46),
    @('I43', 87),
    @('I48', 105),
    @('I49', 67),
    @('I51', 85),
    @('I52', 197),
    @('I53', 109),
    @('I54', 215),
    @('I60', 49),
    @('H63', 197),
    @('I63', 40),
    @('I64', 88),
    @('I67', 366),
    @('I73', 79),
    @('I75', 30),
    @('I79', 246),
    @('I83', 192),
    @('I84', 76),
    @('I85', 445),
    @('I86', 58),
    @('I88', 81),
    @('I89', 104),
    @('I90', 109),
    @('I92', 29),
    @('I93', 58),
    @('I94', 86),
    @('I95', 159),
    @('I96', 117),
    @('I97', 73),
    @('I99', 174),
    @('H101', 25972),
    @('I101', 9622)
)

$updates['South Chicago'] = @(
    @('I2', 69),
    @('I3', 74),
    @('I7', 192)
)

$updates['West Pullman'] = @(
    @('I6', 27),
    @('I7', 159)
)

$updates['Garfield Park'] = @(
    @('I3', 166),
    @('I4', 26),
    @('I7', 452)
)

$updates['Lincoln Park'] = @(
    @('I6', 38),
    @('I7', 67)
)

$updates['Loop'] = @(
    @('I2', 51),
    @('I3', 45),
    @('I6', 104),
    @('I7', 215)
)

$updates['Englewood'] = @(
    @('I2', 199),
    @('I3', 224),
    @('I6', 175),
    @('I7', 642)
)

$updates['Lake View'] = @(
    @('I6', 60),
    @('I7', 105)
)

$updates['South Shore'] = @(
    @('I2', 110),
    @('I3', 183),
    @('I7', 445)
)

$updates['Hermosa'] = @(
    @('I2', 16),
    @('I7', 46)
)

$updates['Douglas'] = @(
    @('I2', 24),
    @('I6', 23)
)

$updates['Chinatown'] = @(
    @('I3', 6),
    @('I7', 57)
)

$updates['Roseland'] = @(
    @('I3', 79),
    @('I7', 246)
)

$updates['Near South Side'] = @(
    @('I2', 22),
    @('I3', 28),
    @('I7', 88)
)

$updates['Chicago Lawn'] = @(
    @('I2', 64),
    @('I7', 241)
)

$updates['West Lawn'] = @(
    @('I2', 15),
    @('I6', 23),
    @('I7', 58)
)

$updates['Little Village'] = @(
    @('I6', 42),
    @('I7', 197)
)

$updates['West Loop'] = @(
    @('I4', 8),
    @('I7', 86)
)

$updates['Brighton Park'] = @(
    @('I2', 40),
    @('I7', 125)
)

$updates['Belmont Cragin'] = @(
    @('I3', 34),
    @('I7', 160)
)

$updates['Portage Park'] = @(
    @('I2', 25),
    @('I7', 79)
)

$updates['Albany Park'] = @(
    @('I3', 35),
    @('I7', 93)
)

$updates['West Town'] = @(
    @('I4', 4),
    @('I7', 73)
)

$updates['West Elsdon'] = @(
    @('I3', 5),
    @('I7', 29)
)

$updates['United Center'] = @(
    @('I6', 31),
    @('I7', 81)
)

$updates['Austin'] = @(
    @('I3', 169),
    @('I4', 39),
    @('I6', 194),
    @('I7', 614)
)

$updates['Streeterville'] = @(
    @('I2', 9),
    @('I4', 36),
    @('I7', 58)
)

$updates['Pullman'] = @(
    @('I2', 11),
    @('I7', 30)
)

$updates['Washington Heights'] = @(
    @('I3', 19),
    @('I7', 109)
)

$updates['Little Italy, UIC'] = @(
    @('I6', 41),
    @('I7', 85)
)

$updates['Morgan Park'] = @(
    @('I6', 19),
    @('I7', 49)
)

$updates['Hyde Park'] = @(
    @('I6', 50),
    @('I7', 87)
)

$updates['Logan Square'] = @(
    @('I2', 19),
    @('I3', 31),
    @('I5', 3),
    @('I6', 48),
    @('I7', 109)
)

$updates['Auburn Gresham'] = @(
    @('I2', 109),
    @('I3', 97),
    @('I6', 81),
    @('I7', 320)
)

$updates['Archer Heights'] = @(
    @('I2', 13),
    @('I7', 39)
)

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $updates[$sheetName]) {
        $cellRef = $pair[0]
        $newValue = $pair[1]
        $ws.Range($cellRef).Value = $newValue
    }
}

$wb.Save()